$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45987
$ws.Cells.Item(2, 2).Value = 11756.64961601
$ws.Cells.Item(2, 3).Value = 10724.0625833776
$ws.Cells.Item(2, 4).Value = 17163.4
$ws.Cells.Item(2, 5).Value = 7176.55096159328
$ws.Cells.Item(2, 6).Value = 30.7172310404546

$ws.Cells.Item(3, 1).Value = 45988
$ws.Cells.Item(3, 2).Value = 11697.8072373275
$ws.Cells.Item(3, 3).Value = 10978.5634005313
$ws.Cells.Item(3, 4).Value = 9003.4
$ws.Cells.Item(3, 5).Value = 6980.21044280814
$ws.Cells.Item(3, 6).Value = 373.140576805809

$ws.Cells.Item(4, 1).Value = 45989
$ws.Cells.Item(4, 2).Value = 11606.1793402161
$ws.Cells.Item(4, 3).Value = 10386.2802690438
$ws.Cells.Item(4, 4).Value = 9003.4
$ws.Cells.Item(4, 5).Value = 6908.4349440701
$ws.Cells.Item(4, 6).Value = 345.47146721308

$ws.Cells.Item(5, 1).Value = 45990
$ws.Cells.Item(5, 2).Value = 4701.20879247002
$ws.Cells.Item(5, 3).Value = 7287.66405124159
$ws.Cells.Item(5, 4).Value = 9003.4
$ws.Cells.Item(5, 5).Value = 6348.9831583061
$ws.Cells.Item(5, 6).Value = 193.051967064487

$ws.Cells.Item(6, 1).Value = 45991
$ws.Cells.Item(6, 2).Value = 4279.73273919308
$ws.Cells.Item(6, 3).Value = 6864.62234033915
$ws.Cells.Item(6, 4).Value = 9003.4
$ws.Cells.Item(6, 5).Value = 6017.8530657882
$ws.Cells.Item(6, 6).Value = 161.628141921973

$ws.Cells.Item(7, 1).Value = 45992
$ws.Cells.Item(7, 2).Value = 11775.3234136769
$ws.Cells.Item(7, 3).Value = 10377.9739801847
$ws.Cells.Item(7, 4).Value = 11040.26
$ws.Cells.Item(7, 5).Value = 7488.53560131314
$ws.Cells.Item(7, 6).Value = 284.427065895745

$ws.Cells.Item(8, 1).Value = 45993
$ws.Cells.Item(8, 2).Value = 11775.3234136769
$ws.Cells.Item(8, 3).Value = 10474.3728876477
$ws.Cells.Item(8, 4).Value = 11040.26
$ws.Cells.Item(8, 5).Value = 7488.53560131314
$ws.Cells.Item(8, 6).Value = 288.443687040036

$ws.Cells.Item(9, 1).Value = 45994
$ws.Cells.Item(9, 2).Value = 11775.3234136769
$ws.Cells.Item(9, 3).Value = 10530.9967483467
$ws.Cells.Item(9, 4).Value = 11040.26
$ws.Cells.Item(9, 5).Value = 7488.53560131314
$ws.Cells.Item(9, 6).Value = 290.803014569161

$ws.Cells.Item(10, 1).Value = 45995
$ws.Cells.Item(10, 2).Value = 11775.3234136769
$ws.Cells.Item(10, 3).Value = 10429.9470237728
$ws.Cells.Item(10, 4).Value = 11040.26
$ws.Cells.Item(10, 5).Value = 7488.53560131314
$ws.Cells.Item(10, 6).Value = 286.592609378582

$ws.Cells.Item(11, 1).Value = 45996
$ws.Cells.Item(11, 2).Value = 11775.3234136769
$ws.Cells.Item(11, 3).Value = 9758.75733173775
$ws.Cells.Item(11, 4).Value = 11040.26
$ws.Cells.Item(11, 5).Value = 7488.53560131314
$ws.Cells.Item(11, 6).Value = 258.626372210454

$ws.Cells.Item(12, 1).Value = 45997
$ws.Cells.Item(12, 2).Value = 4605.75571869631
$ws.Cells.Item(12, 3).Value = 6441.40400727784
$ws.Cells.Item(12, 4).Value = 11040.26
$ws.Cells.Item(12, 5).Value = 7098.35841878054
$ws.Cells.Item(12, 6).Value = 104.145934419099

$ws.Cells.Item(13, 1).Value = 45998
$ws.Cells.Item(13, 2).Value = 4450.60726673227
$ws.Cells.Item(13, 3).Value = 7003.83089365807
$ws.Cells.Item(13, 4).Value = 11040.26
$ws.Cells.Item(13, 5).Value = 7083.29109216126
$ws.Cells.Item(13, 6).Value = 126.952582742472

$ws.Cells.Item(14, 1).Value = 45999
$ws.Cells.Item(14, 2).Value = 4752.2157406543
$ws.Cells.Item(14, 3).Value = 7361.64556246767
$ws.Cells.Item(14, 4).Value = 11040.26
$ws.Cells.Item(14, 5).Value = 7517.01398298563
$ws.Cells.Item(14, 6).Value = 159.933314393887

$ws.Cells.Item(15, 1).Value = 46000
$ws.Cells.Item(15, 2).Value = 12413.0231244078
$ws.Cells.Item(15, 3).Value = 11298.4245867365
$ws.Cells.Item(15, 4).Value = 11040.26
$ws.Cells.Item(15, 5).Value = 7940.92664178886
$ws.Cells.Item(15, 6).Value = 341.628801188558
